# Points per game workbook update:
#  - Add "Super Bowl Winners" and "Sheet2" worksheets
#  - Populate Super Bowl Winners with SB/Date/Winner/Score/Offense/Defense table
#  - Apply custom date number format to the Date column
#  - Restore original active-sheet/window selection state

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets at the end of the workbook, in order.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sbSheet = $wb.Worksheets.Add($null, $lastSheet)
$sbSheet.Name = "Super Bowl Winners"

$sheet2 = $wb.Worksheets.Add($null, $sbSheet)
$sheet2.Name = "Sheet2"

# ---------------------------------------------------------------------------
# 2. Populate the "Super Bowl Winners" sheet.
# ---------------------------------------------------------------------------
$headers = @("SB", "Date", "Winner", "Score", "Offense", "Defense")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $sbSheet.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$rows = @(
    @("XXXVII",   37647, "Tampa Bay Buccaneers", "48–21 vs Oakland Raiders",        "18th", "1st"),
    @("XXXVIII",  38018, "New England Patriots", "32–29 vs Carolina Panthers",      "12th", "1st"),
    @("XXXIX",    38389, "New England Patriots", "24–21 vs Philadelphia Eagles",    "4th",  "2nd"),
    @("XL",       38753, "Pittsburgh Steelers",  "21–10 vs Seattle Seahawks",       "9th",  "3rd"),
    @("XLI",      39117, "Indianapolis Colts",   "29–17 vs Chicago Bears",          "2nd",  "23rd"),
    @("XLII",     39481, "New York Giants",      "17–14 vs New England Patriots",   "14th", "17th"),
    @("XLIII",    39845, "Pittsburgh Steelers",  "27–23 vsArizona Cardinals",       "20th", "1st"),
    @("XLIV",     40216, "New Orleans Saints",   "31–17 vs Indianapolis Colts",     "1st",  "20th"),
    @("XLV",      40580, "Green Bay Packers",    "31–25 vs Pittsburgh Steelers",    "10th", "2nd"),
    @("XLVI",     40944, "New York Giants",      "21–17 vs New England Patriots",   "9th",  "25th"),
    @("XLVII",    41308, "?",                    "at New Orleans, Louisiana",       $null,  $null)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowNum = $r + 2
    $row = $rows[$r]
    $sbSheet.Cells.Item($rowNum, 1).Value = $row[0]
    $sbSheet.Cells.Item($rowNum, 2).Value = $row[1]
    $sbSheet.Cells.Item($rowNum, 3).Value = $row[2]
    $sbSheet.Cells.Item($rowNum, 4).Value = $row[3]
    if ($row[4]) { $sbSheet.Cells.Item($rowNum, 5).Value = $row[4] }
    if ($row[5]) { $sbSheet.Cells.Item($rowNum, 6).Value = $row[5] }
}

# Date column formatting (adds numFmtId 164 / cellXf style index 5)
$sbSheet.Range("B1:B12").NumberFormat = "[$-409]mmmm\ d\,\ yyyy;@"

# Column widths to match the authored layout
$sbSheet.Columns.Item(2).ColumnWidth = 15 - 0.8333333333333333
$sbSheet.Columns.Item(3).ColumnWidth = 20 - 0.8333333333333333
$sbSheet.Columns.Item(4).ColumnWidth = 26.33203125 - 0.8333333333333333

# Match the page margins used throughout the rest of the workbook (inches)
foreach ($sheet in @($sbSheet, $sheet2)) {
    $sheet.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
    $sheet.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
    $sheet.PageSetup.TopMargin = $excel.InchesToPoints(1)
    $sheet.PageSetup.BottomMargin = $excel.InchesToPoints(1)
    $sheet.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
    $sheet.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)
}

# Hide the ruler (matches showRuler="0" used on every sheet in this workbook)
# on both of the newly-created sheets.
$sbSheet.Activate()
$sbSheet.Range("A1:F12").Select() | Out-Null
$excel.ActiveWindow.DisplayRuler = $false

$sheet2.Activate()
$excel.ActiveWindow.DisplayRuler = $false

# ---------------------------------------------------------------------------
# 3. Restore the original active sheet/window state (Offense stays active).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
